$wb = $excel.ActiveWorkbook

$wsTypography = $wb.Worksheets.Item("Typography")
$wsTranslation = $wb.Worksheets.Item("Translation")

# Typography sheet: set the Wildcard Characters for the LCD_Default font (row 7, column G)
# The value looks numeric, so force text formatting first to avoid Excel
# silently converting "0123456789." into the number 123456789, then restore
# the default (unstyled) look of the cell.
$wsTypography.Range("G7").NumberFormat = "@"
$wsTypography.Range("G7").Value = "0123456789."
$wsTypography.Range("G7").Style = "Normal"

# Translation sheet: update existing rows
$wsTranslation.Range("F5").Value = "Voltage2"
$wsTranslation.Range("F6").Value = "<value>"

# Translation sheet: add new rows for voltage/current readouts
$wsTranslation.Range("B7").Value = "voltage"
$wsTranslation.Range("C7").Value = "Small"
$wsTranslation.Range("D7").Value = "Left"
$wsTranslation.Range("E7").Value = "LTR"
$wsTranslation.Range("F7").Value = "Voltage"

$wsTranslation.Range("B8").Value = "current"
$wsTranslation.Range("C8").Value = "Small"
$wsTranslation.Range("D8").Value = "Left"
$wsTranslation.Range("E8").Value = "LTR"
$wsTranslation.Range("F8").Value = "Current"

$wsTranslation.Range("B9").Value = "voltUnit"
$wsTranslation.Range("C9").Value = "Default"
$wsTranslation.Range("D9").Value = "Left"
$wsTranslation.Range("E9").Value = "LTR"
$wsTranslation.Range("F9").Value = "V"

$wsTranslation.Range("B10").Value = "ampUnit"
$wsTranslation.Range("C10").Value = "Default"
$wsTranslation.Range("D10").Value = "Left"
$wsTranslation.Range("E10").Value = "LTR"
$wsTranslation.Range("F10").Value = "A"
